{"js": "// Replace each two-digit-multiplication answer cell's text with its\n// updated value. Every data cell in the table holds a unique\n// \"NNxNN=NNNN\" string, so an exact, case-sensitive whole-text search\n// for the OLD string unambiguously finds the ONE cell/run to update.\nconst replacements = [\n  [\"59\u00d724=1416\", \"50\u00d749=2450\"],\n  [\"43\u00d796=4128\", \"96\u00d768=6528\"],\n  [\"62\u00d769=4278\", \"30\u00d798=2940\"],\n  [\"52\u00d773=3796\", \"61\u00d785=5185\"],\n  [\"28\u00d773=2044\", \"50\u00d712=600\"],\n  [\"37\u00d783=3071\", \"36\u00d757=2052\"],\n  [\"56\u00d783=4648\", \"66\u00d790=5940\"],\n  [\"81\u00d747=3807\", \"74\u00d796=7104\"],\n  [\"17\u00d757=969\", \"16\u00d770=1120\"],\n  [\"19\u00d799=1881\", \"39\u00d736=1404\"],\n  [\"89\u00d753=4717\", \"88\u00d754=4752\"],\n  [\"87\u00d759=5133\", \"50\u00d743=2150\"],\n  [\"77\u00d745=3465\", \"37\u00d738=1406\"],\n  [\"64\u00d758=3712\", \"76\u00d761=4636\"],\n  [\"34\u00d764=2176\", \"31\u00d796=2976\"],\n  [\"74\u00d762=4588\", \"73\u00d712=876\"],\n  [\"63\u00d780=5040\", \"95\u00d789=8455\"],\n  [\"84\u00d730=2520\", \"28\u00d790=2520\"],\n  [\"31\u00d737=1147\", \"52\u00d774=3848\"],\n  [\"99\u00d717=1683\", \"97\u00d765=6305\"],\n  [\"81\u00d756=4536\", \"62\u00d731=1922\"],\n  [\"64\u00d775=4800\", \"64\u00d787=5568\"],\n  [\"17\u00d726=442\", \"80\u00d742=3360\"],\n  [\"50\u00d798=4900\", \"34\u00d716=544\"],\n  [\"87\u00d756=4872\", \"96\u00d769=6624\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Could not find text: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace each two-digit-multiplication answer cell's text with its\n# updated value. Every data cell in the table holds a unique\n# \"NNxNN=NNNN\" string, so an exact whole-text Find/Replace on the OLD\n# string unambiguously targets the ONE cell that needs updating.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"59\u00d724=1416\", \"50\u00d749=2450\"),\n    @(\"43\u00d796=4128\", \"96\u00d768=6528\"),\n    @(\"62\u00d769=4278\", \"30\u00d798=2940\"),\n    @(\"52\u00d773=3796\", \"61\u00d785=5185\"),\n    @(\"28\u00d773=2044\", \"50\u00d712=600\"),\n    @(\"37\u00d783=3071\", \"36\u00d757=2052\"),\n    @(\"56\u00d783=4648\", \"66\u00d790=5940\"),\n    @(\"81\u00d747=3807\", \"74\u00d796=7104\"),\n    @(\"17\u00d757=969\", \"16\u00d770=1120\"),\n    @(\"19\u00d799=1881\", \"39\u00d736=1404\"),\n    @(\"89\u00d753=4717\", \"88\u00d754=4752\"),\n    @(\"87\u00d759=5133\", \"50\u00d743=2150\"),\n    @(\"77\u00d745=3465\", \"37\u00d738=1406\"),\n    @(\"64\u00d758=3712\", \"76\u00d761=4636\"),\n    @(\"34\u00d764=2176\", \"31\u00d796=2976\"),\n    @(\"74\u00d762=4588\", \"73\u00d712=876\"),\n    @(\"63\u00d780=5040\", \"95\u00d789=8455\"),\n    @(\"84\u00d730=2520\", \"28\u00d790=2520\"),\n    @(\"31\u00d737=1147\", \"52\u00d774=3848\"),\n    @(\"99\u00d717=1683\", \"97\u00d765=6305\"),\n    @(\"81\u00d756=4536\", \"62\u00d731=1922\"),\n    @(\"64\u00d775=4800\", \"64\u00d787=5568\"),\n    @(\"17\u00d726=442\", \"80\u00d742=3360\"),\n    @(\"50\u00d798=4900\", \"34\u00d716=544\"),\n    @(\"87\u00d756=4872\", \"96\u00d769=6624\"),\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1          # wdFindContinue\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n\n    $found = $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n    if (-not $found) {\n        throw \"Could not find text: $oldText\"\n    }\n}\n"}
